$wb = $excel.ActiveWorkbook

# --- Sheet "Sydney": remove the old first data row (Atlantic Dawn / row 2)
# which shifts every later arrival up by one row, and refresh the
# "Last Updated" timestamp (column H) for all remaining data rows.
$wsSydney = $wb.Worksheets.Item("Sydney")
$wsSydney.Rows.Item(2).Delete()

$lastRowSydney = $wsSydney.Cells.Item(1, 8).End(-4121).Row  # xlDown = -4121
for ($r = 2; $r -le $lastRowSydney; $r++) {
    $wsSydney.Cells.Item($r, 8).Value = "2026-02-18 01:02"
}

# --- Sheet "Melbourne": just refresh the "Last Updated" timestamp
# (column G) for every data row.
$wsMelbourne = $wb.Worksheets.Item("Melbourne")
$lastRowMelbourne = $wsMelbourne.Cells.Item(1, 7).End(-4121).Row
for ($r = 2; $r -le $lastRowMelbourne; $r++) {
    $wsMelbourne.Cells.Item($r, 7).Value = "2026-02-18 01:02"
}
